# New crime data collected - weekly CompStat update (66th Precinct)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings) ---
# "Volume 31   Number  12" -> "...13"
$ws.Range("A8").Characters(21, 2).Text = "13"

# "Report Covering the Week  3/18/2024  Through  3/24/2024"
#   -> "...3/25/2024  Through  3/31/2024"
$ws.Range("C9").Characters(27, 9).Text = "3/25/2024"
$ws.Range("C9").Characters(47, 9).Text = "3/31/2024"

# --- Cells that flip between numeric and the text placeholders "0" / "***.*" ---
# Use Copy (value+format) from a donor cell that already carries the desired
# style/type, then overwrite with the new value, so the style index used for
# the re-typed cell matches the rest of the sheet.

# F15: number 1 -> text "0"
$ws.Range("G15").Copy($ws.Range("F15"))

# C17: text "0" -> number 5
$ws.Range("D17").Copy($ws.Range("C17"))
$ws.Range("C17").Value = 5

# D22: text "0" -> number 1
$ws.Range("D18").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1

# E22: text "***.*" -> number -100
$ws.Range("K22").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100

# G22: text "0" -> number 1
$ws.Range("G18").Copy($ws.Range("G22"))
$ws.Range("G22").Value = 1

# H22: text "***.*" -> number -100
$ws.Range("K22").Copy($ws.Range("H22"))
$ws.Range("H22").Value = -100

# F27: number 2 -> text "0"
$ws.Range("G27").Copy($ws.Range("F27"))

# C28: number 2 -> text "0"
$ws.Range("C27").Copy($ws.Range("C28"))

# F33: text "0" -> number 1
$ws.Range("F31").Copy($ws.Range("F33"))
$ws.Range("F33").Value = 1

# I33: text "0" -> number 1
$ws.Range("I28").Copy($ws.Range("I33"))
$ws.Range("I33").Value = 1

# --- Plain numeric value updates ---
$ws.Range("L15").Value = 16.666666666666
$ws.Range("M15").Value = 250

$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = -50
$ws.Range("I16").Value = 22
$ws.Range("J16").Value = 19
$ws.Range("K16").Value = 15.789473684210
$ws.Range("L16").Value = -18.518518518518
$ws.Range("M16").Value = -55.102040816326
$ws.Range("N16").Value = -89

$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 66.666666666666
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = -30
$ws.Range("I17").Value = 48
$ws.Range("J17").Value = 50
$ws.Range("K17").Value = -4
$ws.Range("L17").Value = -7.692307692307
$ws.Range("M17").Value = 50
$ws.Range("N17").Value = -25

$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 133.333333333333
$ws.Range("F18").Value = 13
$ws.Range("H18").Value = 44.444444444444
$ws.Range("I18").Value = 32
$ws.Range("J18").Value = 33
$ws.Range("K18").Value = -3.030303030303
$ws.Range("L18").Value = -28.888888888888
$ws.Range("M18").Value = -64.444444444444
$ws.Range("N18").Value = -92.792792792792

$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -12.5
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = 8.823529411764
$ws.Range("I19").Value = 136
$ws.Range("J19").Value = 138
$ws.Range("K19").Value = -1.449275362318
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 47.826086956521
$ws.Range("N19").Value = -14.465408805031

$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 300
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 27.272727272727
$ws.Range("I20").Value = 44
$ws.Range("J20").Value = 34
$ws.Range("K20").Value = 29.411764705882
$ws.Range("L20").Value = 109.52380952381
$ws.Range("M20").Value = 22.222222222222
$ws.Range("N20").Value = -90.736842105263

$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = 62.5
$ws.Range("F21").Value = 82
$ws.Range("G21").Value = 82
$ws.Range("I21").Value = 291
$ws.Range("J21").Value = 278
$ws.Range("K21").Value = 4.676258992805
$ws.Range("L21").Value = 1.393728222996
$ws.Range("M21").Value = -3.322259136212
$ws.Range("N21").Value = -78.428465530022

$ws.Range("J22").Value = 5
$ws.Range("K22").Value = -60

$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 74
$ws.Range("G24").Value = 87
$ws.Range("H24").Value = -14.942528735632
$ws.Range("I24").Value = 243
$ws.Range("J24").Value = 278
$ws.Range("K24").Value = -12.589928057554
$ws.Range("L24").Value = -8.988764044943
$ws.Range("M24").Value = 16.826923076923

$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -45.454545454545
$ws.Range("G25").Value = 28
$ws.Range("H25").Value = -35.714285714285
$ws.Range("I25").Value = 59
$ws.Range("J25").Value = 85
$ws.Range("K25").Value = -30.588235294117
$ws.Range("L25").Value = -27.160493827160

$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = -14.285714285714
$ws.Range("F26").Value = 29
$ws.Range("G26").Value = 21
$ws.Range("H26").Value = 38.095238095238
$ws.Range("I26").Value = 105
$ws.Range("J26").Value = 70
$ws.Range("K26").Value = 50
$ws.Range("L26").Value = 40
$ws.Range("M26").Value = 1.941747572815

$ws.Range("L27").Value = 25

$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 4
$ws.Range("H28").Value = 33.333333333333
$ws.Range("J28").Value = 9
$ws.Range("K28").Value = 144.444444444444
$ws.Range("L28").Value = 57.142857142857

$ws.Range("L31").Value = -20
